$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 304, pushing the existing rows 304-313
# down to 306-315 (dimension grows from A1:R313 to A1:R315).
$ws.Range("A304:R305").EntireRow.Insert()

# New row 304 data
$ws.Range("A304").Value = 4
$ws.Range("B304").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C304").Value = "Los Lagos"
$ws.Range("D304").Value = 44747
$ws.Range("E304").Value = 10
$ws.Range("F304").Value = 100114014
$ws.Range("G304").Value = "Betarraga"
$ws.Range("H304").Value = "Sin especificar"
$ws.Range("I304").Value = "Primera"
$ws.Range("J304").Value = 220
$ws.Range("K304").Value = 11000
$ws.Range("L304").Value = 11000
$ws.Range("M304").Value = 11000
$ws.Range("N304").Value = "$/malla 15 kilos"
$ws.Range("O304").Value = "Región Metropolitana"
$ws.Range("P304").Value = 733
$ws.Range("Q304").Value = 15
$ws.Range("R304").Value = "Hortaliza"

# New row 305 data
$ws.Range("A305").Value = 4
$ws.Range("B305").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C305").Value = "Los Lagos"
$ws.Range("D305").Value = 44747
$ws.Range("E305").Value = 10
$ws.Range("F305").Value = 100114014
$ws.Range("G305").Value = "Betarraga"
$ws.Range("H305").Value = "Sin especificar"
$ws.Range("I305").Value = "Primera"
$ws.Range("J305").Value = 800
$ws.Range("K305").Value = 1000
$ws.Range("L305").Value = 1200
$ws.Range("M305").Value = 1100
$ws.Range("N305").Value = "$/paquete 5 unidades"
$ws.Range("O305").Value = "Región del Maule"
$ws.Range("P305").Value = 220
$ws.Range("Q305").Value = 5
$ws.Range("R305").Value = "Hortaliza"
